$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-25 Tuesday" "2025-02-26 Wednesday"

Replace-Text "59÷4=14, 3" "16÷8=2, 0"
Replace-Text "51÷5=10, 1" "42÷3=14, 0"
Replace-Text "38÷7=5, 3" "79÷2=39, 1"
Replace-Text "39÷8=4, 7" "48÷2=24, 0"
Replace-Text "84÷9=9, 3" "40÷7=5, 5"

Replace-Text "47÷9=5, 2" "17÷4=4, 1"
Replace-Text "55÷2=27, 1" "78÷3=26, 0"
Replace-Text "64÷5=12, 4" "84÷3=28, 0"
Replace-Text "21÷8=2, 5" "59÷7=8, 3"
Replace-Text "21÷9=2, 3" "23÷7=3, 2"

Replace-Text "38÷2=19, 0" "81÷4=20, 1"
Replace-Text "35÷9=3, 8" "80÷6=13, 2"
Replace-Text "85÷7=12, 1" "65÷3=21, 2"
Replace-Text "31÷9=3, 4" "67÷8=8, 3"
Replace-Text "56÷3=18, 2" "41÷9=4, 5"

Replace-Text "21÷2=10, 1" "95÷2=47, 1"
Replace-Text "83÷8=10, 3" "19÷8=2, 3"
Replace-Text "62÷6=10, 2" "16÷4=4, 0"
Replace-Text "77÷3=25, 2" "79÷3=26, 1"
Replace-Text "43÷3=14, 1" "28÷9=3, 1"

Replace-Text "77÷7=11, 0" "84÷7=12, 0"
Replace-Text "43÷9=4, 7" "41÷8=5, 1"
Replace-Text "65÷8=8, 1" "84÷7=12, 0"
Replace-Text "81÷8=10, 1" "63÷3=21, 0"
Replace-Text "92÷9=10, 2" "95÷9=10, 5"
